# Add a header row ("Ticker" / "weight") above the existing ticker/weight
# table on the active sheet, pushing all existing data (and the trailing
# blank filler rows) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 1; this shifts every existing row (data rows
# 1-30 and the blank filler rows 31-60) down by one, automatically
# extending the used range to row 61.
$ws.Rows("1:1").Insert()

# Populate the new header row.
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Value = "weight"

# Match the saved selection state (cell B1 active).
$ws.Range("B1").Select()
